# Add a "name" column in front of the "id" column on the Analysis sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("Analysis", "Analysis1")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Shift existing columns (A:K) one to the right to make room for the
    # new "name" column, then write the new header in column A.
    $ws.Columns.Item(1).Insert()
    $ws.Cells.Item(1, 1).Value = "name"
}
